$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the old sub-header row (old row 2: (m3/s)/(MW)/(GWh) labels).
# This shifts the data rows (old 3..14) up to become new rows 2..13.
$ws.Rows(2).Delete()

# Rewrite row 1 with the new column headers.
$ws.Cells.Item(1,1).Value = "idx"
$ws.Cells.Item(1,2).Value = "idx2"
$ws.Cells.Item(1,3).Value = "Name"
$ws.Cells.Item(1,4).Value = "Date Start"
$ws.Cells.Item(1,5).Value = "Date End"
$ws.Cells.Item(1,6).Value = "(m3/s)"
$ws.Cells.Item(1,7).Value = "(MW1)"
$ws.Cells.Item(1,8).Value = "(MW2)"
$ws.Cells.Item(1,9).Value = "(GWh) Winter"
$ws.Cells.Item(1,10).Value = "(GWh) Summer"
$ws.Cells.Item(1,11).Value = "(GWh) Year"

# A1:E1 carry no explicit style (default); clear any inherited formatting.
$ws.Range("A1:E1").Font.Size = 10

# F1:K1 use the small 9pt Arial header style (same as the rest of the sheet).
$ws.Range("F1:K1").Font.Size = 9
$ws.Range("F1:K1").Font.Name = "Arial"

$ws.Range("A2:K2").Select()
